$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item('Overview')
# Remove all existing hyperlinks up front; Insert() does not shift
# hyperlink anchors the way it shifts cell data, so we rebuild them
# from scratch after the row insert below.
$ws.Hyperlinks.Delete()

# Insert a new row 9 (pushes old row 9 '.localization-config' down to row 10)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value2 = 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md'
$ws.Range("B9").Value2 = 'Ready for handoff'
$ws.Range("C9").Value2 = 'Ready for handoff'

# Re-create hyperlinks: rows 2..8 (and extra cols) unchanged, plus new row 9
# source-file link, plus row 9 target-xlf link (zh-cn/de-de only), plus the
# '.localization-config' link which is now on row 10.
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/e2e/2cf53653-fad0-4d68-9bb1-208c798ec3ad.md', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/7c3ff5a4792aea27063cc1ec5c2bd5da775ae051/e2e/20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md', "", "", '20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/fb5bdc908963c21d8f58d0a2f068e664c4bf375a/e2e/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/831315f2-da6a-4fc4-b117-60389f3c6074.md', "", "", '831315f2-da6a-4fc4-b117-60389f3c6074.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md', "", "", 'ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/32067f598ff8e2e6f57aa0a46fea662898476080/e2e/0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md', "", "", '0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/194ec55d33c94f4d591d75bbcd2c5146b4ebc0e9/e2e/984aa78d-5405-411a-8b9f-4541e5a8b93a.md', "", "", '984aa78d-5405-411a-8b9f-4541e5a8b93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/.localization-config', "", "", '.localization-config') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/0a8eed458612eda148562c46b5ec8fb9fc84407f/e2e/f1426313-4c95-4d89-9ab0-d616ca86b43c.md', "", "", 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md') | Out-Null

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item('zh-cn')
# Remove all existing hyperlinks up front; Insert() does not shift
# hyperlink anchors the way it shifts cell data, so we rebuild them
# from scratch after the row insert below.
$ws.Hyperlinks.Delete()

# Insert a new row 9 (pushes old row 9 '.localization-config' down to row 10)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value2 = 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md'
$ws.Range("B9").Value2 = 'Ready for handoff'
$ws.Range("C9").Value2 = 'f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.zh-cn.xlf'
$ws.Range("D9").Value2 = '2016-02-26 04:59:30'
$ws.Range("G9").Value2 = '0001-01-01 00:00:00'
$ws.Range("H9").Value2 = 'Include'

# Re-create hyperlinks: rows 2..8 (and extra cols) unchanged, plus new row 9
# source-file link, plus row 9 target-xlf link (zh-cn/de-de only), plus the
# '.localization-config' link which is now on row 10.
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/e2e/2cf53653-fad0-4d68-9bb1-208c798ec3ad.md', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/630f5c620c63521eb08f7d9f80766f101d4dfc89/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.zh-cn.xlf', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e098be5190c5c6499279feabfdab7c813bc3ae93/e2e/2cf53653-fad0-4d68-9bb1-208c798ec3ad.md', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/55cd36a912caa93b79ef811d529de09c47e5919f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.zh-cn.xlf', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/7c3ff5a4792aea27063cc1ec5c2bd5da775ae051/e2e/20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md', "", "", '20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2a8c969f40bb19acfc1e2357b28af3b468f92ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.zh-cn.xlf', "", "", '20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/fb5bdc908963c21d8f58d0a2f068e664c4bf375a/e2e/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d4d7e185e5dc1977a3318dbbbaac2a0aa8dc985/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E4'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c8f23146bbfd696241de54612728556ea849c7e7/e2e/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/712a1cb9324e65806a82deb2f49284339d11b02b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/831315f2-da6a-4fc4-b117-60389f3c6074.md', "", "", '831315f2-da6a-4fc4-b117-60389f3c6074.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c21876685d0bf6742158b34fd7da88917bfaa62/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.zh-cn.xlf', "", "", '831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md', "", "", 'ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C6'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c21876685d0bf6742158b34fd7da88917bfaa62/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.zh-cn.xlf', "", "", 'ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/32067f598ff8e2e6f57aa0a46fea662898476080/e2e/0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md', "", "", '0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C7'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/edde8c6f9fc1626b526ebd9ff410f82f840cc72e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.dcd3d332b0d7f0538e7c30bc8452a132e6c349a4.zh-cn.xlf', "", "", '0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.dcd3d332b0d7f0538e7c30bc8452a132e6c349a4.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/194ec55d33c94f4d591d75bbcd2c5146b4ebc0e9/e2e/984aa78d-5405-411a-8b9f-4541e5a8b93a.md', "", "", '984aa78d-5405-411a-8b9f-4541e5a8b93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C8'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0168fd15f0e72885ec0025220834e3971d9eb5a0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.zh-cn.xlf', "", "", '984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/.localization-config', "", "", '.localization-config') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/0a8eed458612eda148562c46b5ec8fb9fc84407f/e2e/f1426313-4c95-4d89-9ab0-d616ca86b43c.md', "", "", 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C9'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a8eed458612eda148562c46b5ec8fb9fc84407f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.zh-cn.xlf', "", "", 'f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.zh-cn.xlf') | Out-Null

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item('de-de')
# Remove all existing hyperlinks up front; Insert() does not shift
# hyperlink anchors the way it shifts cell data, so we rebuild them
# from scratch after the row insert below.
$ws.Hyperlinks.Delete()

# Insert a new row 9 (pushes old row 9 '.localization-config' down to row 10)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value2 = 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md'
$ws.Range("B9").Value2 = 'Ready for handoff'
$ws.Range("C9").Value2 = 'f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.de-de.xlf'
$ws.Range("D9").Value2 = '2016-02-26 04:59:40'
$ws.Range("G9").Value2 = '0001-01-01 00:00:00'
$ws.Range("H9").Value2 = 'Include'

# Re-create hyperlinks: rows 2..8 (and extra cols) unchanged, plus new row 9
# source-file link, plus row 9 target-xlf link (zh-cn/de-de only), plus the
# '.localization-config' link which is now on row 10.
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/e2e/2cf53653-fad0-4d68-9bb1-208c798ec3ad.md', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4a61a14fff073195361f3be1af8c3352f2528d1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.de-de.xlf', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8f3e8980a3f06c205d28e64ba15a36d292c12bbf/e2e/2cf53653-fad0-4d68-9bb1-208c798ec3ad.md', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/362c99c35dfa29b530a2e2e1a6503439a7d2549d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.de-de.xlf', "", "", '2cf53653-fad0-4d68-9bb1-208c798ec3ad.852dda83b75cf2422d84d676ac708a770c1be604.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/7c3ff5a4792aea27063cc1ec5c2bd5da775ae051/e2e/20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md', "", "", '20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d526f8aeac246ce3580b2cd9a1ae89a61b1a3d39/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.de-de.xlf', "", "", '20bfff53-8ddc-4c6b-841a-7d0f35d2ba6d.550666d79c2cc553a731e154650f5b8682684363.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/fb5bdc908963c21d8f58d0a2f068e664c4bf375a/e2e/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/39e869b965ea149e9cffd8618ee804cbe692d191/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E4'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d99a50cf9f0fc7bf5daf39124b297e3c92bd2f91/e2e/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/7b8965a022667ff80e7fe40292d20d84a6f20c0d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf', "", "", '3f1c9f3e-43f7-49e1-b3c2-b71c60c0e93a.f1c949f73f69c640522426be7525a3a80e789148.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/831315f2-da6a-4fc4-b117-60389f3c6074.md', "", "", '831315f2-da6a-4fc4-b117-60389f3c6074.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a907fa33aeae90be0f13d9ceb0611eb052393f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.de-de.xlf', "", "", '831315f2-da6a-4fc4-b117-60389f3c6074.33096a8b89da196732c33acf66abed619fc16792.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A6'), 'https://github.com/OpenLocalizationTest/oltest/blob/2383b328c80f92809e3a02927ffa26cd4f4610c3/e2e/ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md', "", "", 'ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C6'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a907fa33aeae90be0f13d9ceb0611eb052393f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.de-de.xlf', "", "", 'ec9ae7f3-106a-4b1a-bd14-737fe3dec9c4.fc1d393ccc271987bd0fbbef8d59cb1dfbb601f9.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A7'), 'https://github.com/OpenLocalizationTest/oltest/blob/32067f598ff8e2e6f57aa0a46fea662898476080/e2e/0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md', "", "", '0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C7'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a7a6260ce375caf857e54c4d82513e9735de216e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.dcd3d332b0d7f0538e7c30bc8452a132e6c349a4.de-de.xlf', "", "", '0cd4e0fa-534a-4c0b-bf1d-b9b57b1db9e7.dcd3d332b0d7f0538e7c30bc8452a132e6c349a4.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A8'), 'https://github.com/OpenLocalizationTest/oltest/blob/194ec55d33c94f4d591d75bbcd2c5146b4ebc0e9/e2e/984aa78d-5405-411a-8b9f-4541e5a8b93a.md', "", "", '984aa78d-5405-411a-8b9f-4541e5a8b93a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C8'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c53ed2ba9ea5d1836ddbc30dfd5aa29f71aeb127/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.de-de.xlf', "", "", '984aa78d-5405-411a-8b9f-4541e5a8b93a.38149897969e365477930dd55e78ffc1469ed844.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/62ff95d134f7a5ad91f5129fdb62b14914703f49/.localization-config', "", "", '.localization-config') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A9'), 'https://github.com/OpenLocalizationTest/oltest/blob/0a8eed458612eda148562c46b5ec8fb9fc84407f/e2e/f1426313-4c95-4d89-9ab0-d616ca86b43c.md', "", "", 'f1426313-4c95-4d89-9ab0-d616ca86b43c.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C9'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a8eed458612eda148562c46b5ec8fb9fc84407f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.de-de.xlf', "", "", 'f1426313-4c95-4d89-9ab0-d616ca86b43c.0a8eed458612eda148562c46b5ec8fb9fc84407f.de-de.xlf') | Out-Null

Write-Host "Done generating handback report."
